# Apply "Added BAARD stage 2, line, watermark to CIFAR10" edit:
#  - update a handful of existing Acc_on_adv(lid)/FPR(lid) values (col F / L)
#  - append two new attack groups (rows 22-24 "line", rows 25-27 "watermark")
#  - merge the new attack-name column cells, matching the existing groups

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- existing-cell value corrections -------------------------------------
$ws.Range("F5").Value  = 72
$ws.Range("F6").Value  = 85.39999999999998
$ws.Range("L6").Value  = 0
$ws.Range("F7").Value  = 95.40000000000001
$ws.Range("F8").Value  = 95.8
$ws.Range("F11").Value = 25.6
$ws.Range("L11").Value = 2.4
$ws.Range("F12").Value = 48.2
$ws.Range("L12").Value = 0.8
$ws.Range("F13").Value = 68.7
$ws.Range("L13").Value = 7.5
$ws.Range("F14").Value = 44.8
$ws.Range("L14").Value = 10.3
$ws.Range("F16").Value = 4.100000000000001
$ws.Range("L16").Value = 0.1
$ws.Range("F17").Value = 15.7
$ws.Range("L17").Value = 1.2
$ws.Range("F21").Value = 7.000000000000001

# ---- new row 22: "line", Adv_param 0 --------------------------------------
$ws.Range("A22").Value = "line"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 81.3
$ws.Range("D22").Value = 81.69999999999997
$ws.Range("E22").Value = 96.90000000000001
$ws.Range("F22").Value = 73.5
$ws.Range("G22").Value = 81.39999999999998
$ws.Range("H22").Value = 74.09999999999999
$ws.Range("I22").Value = 3.5
$ws.Range("J22").Value = 6.4
$ws.Range("K22").Value = 5.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0.4
$ws.Range("N22").Value = 0.4

# ---- new row 23: Adv_param 0.5 --------------------------------------------
$ws.Range("B23").Value = 0.5
$ws.Range("C23").Value = 95
$ws.Range("D23").Value = 95.09999999999999
$ws.Range("E23").Value = 96.90000000000001
$ws.Range("F23").Value = 94.7
$ws.Range("G23").Value = 94.90000000000001
$ws.Range("H23").Value = 94.7
$ws.Range("I23").Value = 3.9
$ws.Range("J23").Value = 6.4
$ws.Range("K23").Value = 5.2
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0.4
$ws.Range("N23").Value = 0.4

# ---- new row 24: Adv_param 1 ----------------------------------------------
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 83
$ws.Range("D24").Value = 82.3
$ws.Range("E24").Value = 82.3
$ws.Range("F24").Value = 76.5
$ws.Range("G24").Value = 82.59999999999999
$ws.Range("H24").Value = 77.10000000000002
$ws.Range("I24").Value = 4.2
$ws.Range("J24").Value = 6
$ws.Range("K24").Value = 5.2
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0.4
$ws.Range("N24").Value = 0.3

# ---- new row 25: "watermark", Adv_param 0.3 --------------------------------
$ws.Range("A25").Value = "watermark"
$ws.Range("B25").Value = 0.3
$ws.Range("C25").Value = 91.7
$ws.Range("D25").Value = 91.90000000000001
$ws.Range("E25").Value = 94.59999999999999
$ws.Range("F25").Value = 91.40000000000001
$ws.Range("G25").Value = 92.5
$ws.Range("H25").Value = 91.40000000000001
$ws.Range("I25").Value = 4.3
$ws.Range("J25").Value = 6.600000000000001
$ws.Range("K25").Value = 5.2
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0.4
$ws.Range("N25").Value = 0.3

# ---- new row 26: Adv_param 0.6 ---------------------------------------------
$ws.Range("B26").Value = 0.6
$ws.Range("C26").Value = 40.40000000000001
$ws.Range("D26").Value = 41.6
$ws.Range("E26").Value = 51.2
$ws.Range("F26").Value = 40.3
$ws.Range("G26").Value = 41.4
$ws.Range("H26").Value = 40.3
$ws.Range("I26").Value = 3.9
$ws.Range("J26").Value = 6.1
$ws.Range("K26").Value = 5.2
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0.4
$ws.Range("N26").Value = 0.4

# ---- new row 27: Adv_param 1 -----------------------------------------------
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 0.6
$ws.Range("D27").Value = 2.2
$ws.Range("E27").Value = 0.4
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0.1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 3.8
$ws.Range("J27").Value = 6.2
$ws.Range("K27").Value = 5.2
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0.4
$ws.Range("N27").Value = 0.3

# ---- carry over the bold/bordered/centered header style (style index used
#      by the existing A/B "group" columns, e.g. A17:B17) onto the new cells
$ws.Range("A17:B17").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)
$ws.Range("A25:B25").PasteSpecial(-4122)

$ws.Range("A18:B18").Copy()
$ws.Range("A23:B23").PasteSpecial(-4122)
$ws.Range("A24:B24").PasteSpecial(-4122)
$ws.Range("A26:B26").PasteSpecial(-4122)
$ws.Range("A27:B27").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- merge the attack-name column for the two new groups -------------------
$ws.Range("A22:A24").Merge()
$ws.Range("A25:A27").Merge()

# merging redraws per-cell borders on the merged range (outer box only), which
# splits the single shared style used across the sheet into extra variants;
# re-stamp the original uniform style so the new rows match the existing
# merged groups (A4:A8, A9:A12, ...) byte-for-byte in style usage
$ws.Range("A17:B17").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)
$ws.Range("A25:B25").PasteSpecial(-4122)

$ws.Range("A18:B18").Copy()
$ws.Range("A23:B23").PasteSpecial(-4122)
$ws.Range("A24:B24").PasteSpecial(-4122)
$ws.Range("A26:B26").PasteSpecial(-4122)
$ws.Range("A27:B27").PasteSpecial(-4122)

$excel.CutCopyMode = $false
